$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.270540475845337
$ws.Range("B1").Value = 2.485193729400635
$ws.Range("C1").Value = 2.48551607131958
$ws.Range("D1").Value = 2.142321825027466
$ws.Range("E1").Value = 1.537627935409546
